$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Programming Languages:" skills line.
#    " Java, C++/C, Python, " -> " C++, C, Java, Python, "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Java, C++/C, Python,", $true, $false, $false, $false, $false,
    $true, 1, $false, "C++, C, Java, Python,", 2) | Out-Null

#    "SQL, HTML, JavaScript, Processing" -> "HTML, JavaScript, SQL"
$d.Content.Find.Execute(
    "SQL, HTML, JavaScript, Processing", $true, $false, $false, $false, $false,
    $true, 1, $false, "HTML, JavaScript, SQL", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a new EXPERIENCE bullet ("Fall 2024 ... virtual memory system")
#    right before the existing "June 2024" bullet.
# ---------------------------------------------------------------------------
$juneFound = $d.Content.Find.Execute("June 2024" + [char]9 + "built", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("June 2024")) {
        $targetIndex = $i
        break
    }
}

$junePara = $d.Paragraphs($targetIndex)
$junePara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs($targetIndex)

$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:ind w:left="1260" w:hanging="1260"/>
              <w:jc w:val="left"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="SimSun" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:bCs/>
                <w:color w:val="000000"/>
                <w:kern w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:bidi="ar"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="SimSun" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:bCs/>
                <w:color w:val="000000"/>
                <w:kern w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:bidi="ar"/>
              </w:rPr>
              <w:t>Fall 2024</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="SimSun" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:bCs/>
                <w:color w:val="000000"/>
                <w:kern w:val="0"/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
                <w:lang w:bidi="ar"/>
              </w:rPr>
              <w:tab/>
              <w:t>Designed and implemented a simulated virtual memory system (TLB/Cache/Mem) in C</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 3. Capitalise "built " -> "Built " (both EXPERIENCE bullets that start a
#    sentence with it).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "built ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Built ", 2) | Out-Null
$d.Content.Find.Execute(
    "built ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Built ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. HOBBIES line.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Ultimate Frisbee, Swimming, Football, Ping-pong, Basketball, Badminton, Black Myth: Wukong, NFL Madden",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Ultimate Frisbee, Swimming, football, Pingpong, Basketball, Badminton, Madden", 2) | Out-Null
